$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.26"
$ws.Range("E2").Value = "'5.39%"

$ws.Range("D3").Value = "'34.81"
$ws.Range("E3").Value = "'12.28%"

$ws.Range("D4").Value = "'5.175"
$ws.Range("E4").Value = "'4.35%"

$ws.Range("D5").Value = "'0.07799"
$ws.Range("E5").Value = "'6.11%"

$ws.Range("D6").Value = "'2.306"
$ws.Range("E6").Value = "'-0.43%"

$ws.Range("D7").Value = "'8.056"
$ws.Range("E7").Value = "'4.38%"

$ws.Range("D8").Value = "'3.991"
$ws.Range("E8").Value = "'6.71%"

$ws.Range("D9").Value = "'0.9247"
$ws.Range("E9").Value = "'1.41%"

$ws.Range("E10").Value = "'9.29%"

$ws.Range("D11").Value = "'0.1827"
$ws.Range("E11").Value = "'7.54%"

$ws.Range("D12").Value = "'0.08519"
$ws.Range("E12").Value = "'3.11%"

$ws.Range("D13").Value = "'0.03390"
$ws.Range("E13").Value = "'9.05%"

$ws.Range("D14").Value = "'0.09914"
$ws.Range("E14").Value = "'-0.69%"

$ws.Range("D15").Value = "'0.001494"
$ws.Range("E15").Value = "'-0.14%"

$ws.Range("D16").Value = "'0.04654"
$ws.Range("E16").Value = "'2.95%"

$ws.Range("D17").Value = "'0.005807"
$ws.Range("E17").Value = "'1.17%"

$ws.Range("D18").Value = "'3.472"
$ws.Range("E18").Value = "'0.02%"

$ws.Range("D19").Value = "'2.103"
$ws.Range("E19").Value = "'3.17%"

$ws.Range("E20").Value = "'2.87%"

$ws.Range("D21").Value = "'0.1327"
$ws.Range("E21").Value = "'3.04%"

$ws.Range("D22").Value = "'4.553"
$ws.Range("E22").Value = "'9.29%"

$ws.Range("D23").Value = "'0.2273"
$ws.Range("E23").Value = "'7.01%"

$ws.Range("D24").Value = "'0.001221"
$ws.Range("E24").Value = "'0.83%"

$ws.Range("D25").Value = "'0.004339"
$ws.Range("E25").Value = "'3.92%"

$ws.Range("E26").Value = "'0.02%"

$ws.Range("D27").Value = "'0.0003403"
$ws.Range("E27").Value = "'0.25%"

$ws.Range("D39").Value = "'0.01745"
$ws.Range("E39").Value = "'11.18%"

$ws.Range("D40").Value = "'0.04745"
$ws.Range("E40").Value = "'6.10%"

$ws.Range("D41").Value = "'0.007684"
$ws.Range("E41").Value = "'4.60%"

$ws.Range("E42").Value = "'6.03%"

$ws.Range("E43").Value = "'-22.39%"

$ws.Range("D44").Value = "'0.002291"
$ws.Range("E44").Value = "'2.25%"

$ws.Range("D45").Value = "'0.009967"
$ws.Range("E45").Value = "'13.78%"

$ws.Range("D46").Value = "'0.00006068"
$ws.Range("E46").Value = "'-0.63%"

$ws.Range("E47").Value = "'0.14%"

$ws.Range("D48").Value = "'5.795"
$ws.Range("E48").Value = "'136.99%"

$ws.Range("D49").Value = "'0.002694"
$ws.Range("E49").Value = "'34.68%"

$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.14%"

$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.14%"
